# #ERM-15 Add all the columns from the import while exporting error data
#
# Fixes sample/placeholder data in the first two data rows of the User
# Registration template and switches the "FromTime"/"To Time" availability
# columns (AC/AD) from Excel time-serial values to plain text values
# (e.g. "9:30 AM" / "10:00 PM") formatted as Text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (HUB Kerry / Mohit Joshi) -------------------------------------
$ws.Range("B3").Value2 = "HUB Kerry"
$ws.Range("F3").Value2 = "Kerry@gmail.com"
$ws.Range("L3").Value2 = "Mohit"
$ws.Range("M3").Value2 = "Joshi"
$ws.Range("N3").Value2 = 9434454873
$ws.Range("S3").Value2 = "Mohit@cdac.com"

# --- Row 4 (UPHC Debbe / Guneet Wagh) ------------------------------------
$ws.Range("B4").Value2 = "UPHC Debbe"
$ws.Range("F4").Value2 = "Debbe@gmail.com"
$ws.Range("L4").Value2 = "Guneet"
$ws.Range("N4").Value2 = 8632362355
$ws.Range("S4").Value2 = "Guneet@gmail.com"

# --- Availability times: store as text "9:30 AM"/"10:00 PM" instead of ---
# --- a numeric time-of-day serial value. --------------------------------
$timeFormat = "@"

$ws.Range("AC3").NumberFormat = $timeFormat
$ws.Range("AC3").Value2 = "9:30 AM"
$ws.Range("AD3").NumberFormat = $timeFormat
$ws.Range("AD3").Value2 = "10:00 PM"

$ws.Range("AC4").NumberFormat = $timeFormat
$ws.Range("AC4").Value2 = "9:30 AM"
$ws.Range("AD4").NumberFormat = $timeFormat
$ws.Range("AD4").Value2 = "10:00 PM"
